$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new Price text in D, new Volume(1h) text in E). $null means D is unchanged.
$updates = @(
    @{ Row = 2; D = "27.945.22"; E = "  +0.81%  " },
    @{ Row = 3; D = "1.880.41"; E = "  +0.26%  " },
    @{ Row = 4; D = "1.017"; E = "  +1.35%  " },
    @{ Row = 5; D = "335.02"; E = "  +0.88%  " },
    @{ Row = 6; D = "1.016"; E = "  +1.23%  " },
    @{ Row = 7; D = "0.4690"; E = "  -0.56%  " },
    @{ Row = 8; D = "0.3916"; E = "  -1.03%  " },
    @{ Row = 9; D = $null; E = "  -1.45%  " },
    @{ Row = 10; D = "0.07959"; E = "  -0.96%  " },
    @{ Row = 11; D = "1.008"; E = "  -1.52%  " },
    @{ Row = 12; D = "21.63"; E = "  -1.03%  " },
    @{ Row = 13; D = "1.884.18"; E = "  -0.47%  " },
    @{ Row = 14; D = "5.958"; E = "  +0.00%  " },
    @{ Row = 15; D = "7.102"; E = "  -0.69%  " },
    @{ Row = 16; D = "1.019"; E = "  +1.28%  " },
    @{ Row = 17; D = "0.06800"; E = "  +2.62%  " },
    @{ Row = 18; D = "87.58"; E = "  +0.42%  " },
    @{ Row = 19; D = "0.00001042"; E = "  -0.60%  " },
    @{ Row = 20; D = "17.04"; E = "  -0.86%  " },
    @{ Row = 21; D = $null; E = "  +1.30%  " },
    @{ Row = 22; D = "27.946.34"; E = "  +0.54%  " },
    @{ Row = 23; D = "5.473"; E = "  -0.57%  " },
    @{ Row = 24; D = "10.98"; E = "  -0.42%  " },
    @{ Row = 25; D = "2.355"; E = "  +2.34%  " },
    @{ Row = 26; D = "2.112.40"; E = "  -0.25%  " },
    @{ Row = 27; D = "159.49"; E = "  +1.99%  " },
    @{ Row = 28; D = "19.90"; E = "  -1.54%  " },
    @{ Row = 29; D = "2.073"; E = "  -1.37%  " },
    @{ Row = 30; D = "5.468"; E = "  -2.03%  " },
    @{ Row = 31; D = "120.91"; E = "  -1.32%  " },
    @{ Row = 32; D = "0.09524"; E = "  -0.32%  " },
    @{ Row = 33; D = "0.9576"; E = "  -1.14%  " },
    @{ Row = 34; D = "3.656"; E = "  +0.67%  " },
    @{ Row = 35; D = "5.330"; E = "  +0.51%  " },
    @{ Row = 36; D = $null; E = "  -7.22%  " },
    @{ Row = 37; D = "0.06119"; E = "  +0.01%  " },
    @{ Row = 38; D = "0.02236"; E = "  -1.20%  " },
    @{ Row = 39; D = "1.203"; E = "  -2.12%  " },
    @{ Row = 40; D = $null; E = "  +1.26%  " },
    @{ Row = 41; D = "8.114"; E = "  -0.85%  " },
    @{ Row = 42; D = "0.5896"; E = "  -1.52%  " },
    @{ Row = 43; D = "0.1895"; E = "  -0.81%  " },
    @{ Row = 44; D = $null; E = "  -0.63%  " },
    @{ Row = 45; D = "1.272"; E = "  +2.01%  " },
    @{ Row = 46; D = "0.5648"; E = "  -1.09%  " },
    @{ Row = 47; D = "12.13"; E = "  -0.72%  " },
    @{ Row = 48; D = "3.399"; E = "  -0.18%  " },
    @{ Row = 49; D = "1.919"; E = "  -0.66%  " },
    @{ Row = 51; D = "113.76"; E = "  +1.34%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        # Force text storage so numeric-looking strings (e.g. "1.017") are not
        # reinterpreted as numbers; ClearFormats afterwards restores the default
        # (unstyled) cell format so no stray style is left behind.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value2 = $u.D
        $ws.Range("D$row").ClearFormats()
    }
    $ws.Range("E$row").Value2 = $u.E
}
